# Fruta / hortaliza, semanal
# Insert a new weekly record before the last row of data, and update the
# (previously) last row with the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 19. This pushes the former row 19 down to row 20,
# leaving the (still unmodified) former row 18 data in place at row 18.
$ws.Rows.Item(19).Insert()

# Copy the former row-18 values (A18:T18) into the newly inserted row 19,
# since row 19 is now a duplicate data point for the week before the new one.
$ws.Range("A19").Value = $ws.Range("A18").Value()
$ws.Range("B19").Value = $ws.Range("B18").Value()
$ws.Range("C19").Value = $ws.Range("C18").Value()
$ws.Range("D19").Value = $ws.Range("D18").Value()
$ws.Range("E19").Value = $ws.Range("E18").Value()
$ws.Range("F19").Value = $ws.Range("F18").Value()
$ws.Range("G19").Value = $ws.Range("G18").Value()
$ws.Range("H19").Value = $ws.Range("H18").Value()
$ws.Range("I19").Value = $ws.Range("I18").Value()
$ws.Range("J19").Value = $ws.Range("J18").Value()
$ws.Range("K19").Value = $ws.Range("K18").Value()
$ws.Range("L19").Value = $ws.Range("L18").Value()
$ws.Range("M19").Value = $ws.Range("M18").Value()
$ws.Range("N19").Value = $ws.Range("N18").Value()
$ws.Range("O19").Value = $ws.Range("O18").Value()
$ws.Range("P19").Value = $ws.Range("P18").Value()
$ws.Range("Q19").Value = $ws.Range("Q18").Value()
$ws.Range("R19").Value = $ws.Range("R18").Value()
$ws.Range("S19").Value = $ws.Range("S18").Value()
$ws.Range("T19").Value = $ws.Range("T18").Value()

# Now overwrite row 18 with the new week's figures (new date, new prices).
$ws.Range("D18").Value = 44876
$ws.Range("N18").Value = 7500
$ws.Range("O18").Value = 8000
$ws.Range("P18").Value = 7750
$ws.Range("S18").Value = 5167
